$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition listing)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 1916
$wsExhibit.Range("F4").Value = 827
$wsExhibit.Range("F5").Value = 843
$wsExhibit.Range("F6").Value = 268

# Sheet "全部类型" (all types combined listing)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1916
$wsAll.Range("F5").Value = 827
$wsAll.Range("F6").Value = 843
$wsAll.Range("F7").Value = 268
